# Rename the inline picture "name" labels (wp:docPr / pic:cNvPr) for the
# three logo images in this document:
#   - Footers(1) [primary footer]   PearsonLogo: image2.png -> image1.png
#   - Footers(2) [first-page footer] PearsonLogo: image2.png -> image1.png
#   - Headers(2) [first-page header] BTec_Logo-Orange: image1.jpg -> image2.jpg
#
# InlineShape has no settable .Name property in the Word object model;
# Shape.Name is what backs <wp:docPr name="...">, so each picture is
# temporarily converted to a floating Shape, renamed, then converted back
# to an inline shape (restoring <wp:inline>) to keep the rest of the
# markup/layout unchanged.

$d = $word.ActiveDocument

function Rename-InlinePicture($range, $newName) {
    $shp = $range.InlineShapes(1)
    $floating = $shp.ConvertToShape()
    $floating.Name = $newName
    [void]$floating.ConvertToInlineShape()
}

$sec = $d.Sections(1)

# Primary footer - PearsonLogo.png (id=1)
Rename-InlinePicture $sec.Footers(1).Range "image1.png"

# First-page footer - PearsonLogo.png (id=2)
Rename-InlinePicture $sec.Footers(2).Range "image1.png"

# First-page header - BTec_Logo-Orange (id=3)
Rename-InlinePicture $sec.Headers(2).Range "image2.jpg"
